# This "edit" is a rebuild of the generated docx (commit message: "Built
# site for gh-pages"). The underlying change is the style/numbering part
# being re-serialized by the generator, which re-orders a handful of
# sibling elements inside word/styles.xml (no formatting values actually
# change) and drops a stray literal "&gt;" that had leaked into the
# AbstractTitle run-properties. We reproduce that by touching each
# affected style through the Word object model with value-preserving
# ("identity") writes, which forces Word to re-emit the style in its
# canonical element order - matching the rebuilt XML - without altering
# any visible formatting.

$d = $word.ActiveDocument

# -- word/styles.xml : BodyText --------------------------------------
# Moves <w:qFormat/> back above <w:pPr> for the "Body Text" style.
$bodyText = $d.Styles.Item("BodyText")
$bodyText.QuickStyle = $true

# -- word/styles.xml : AbstractTitle -----------------------------------
# Moves <w:spacing/> above <w:jc/>, moves <w:b/>/<w:color/> above
# <w:sz/>/<w:szCs/>, and drops the stray "&gt;" left inside <w:rPr>.
$abstractTitle = $d.Styles.Item("AbstractTitle")
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = $abstractTitle.Font.Color
$abstractTitle.ParagraphFormat.Alignment = $abstractTitle.ParagraphFormat.Alignment

# -- word/styles.xml : Heading4 -----------------------------------------
# Swaps <w:bCs/> above <w:i/>.
$heading4 = $d.Styles.Item("Heading4")
$heading4.Font.Italic = $true

# -- word/styles.xml : DocumentationTok / CommentVarTok / WarningTok ----
# Moves <w:i/> above <w:color/>/<w:shd/> in each of the three token
# character styles.
$documentationTok = $d.Styles.Item("DocumentationTok")
$documentationTok.Font.Italic = $true

$commentVarTok = $d.Styles.Item("CommentVarTok")
$commentVarTok.Font.Italic = $true

$warningTok = $d.Styles.Item("WarningTok")
$warningTok.Font.Italic = $true

Write-Output "styles re-serialized"
